$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1150
$ws.Range("I40").Value = 1110.7142
$ws.Range("J40").Value = 1333.3334
$ws.Range("K40").Value = 1110.7142
$ws.Range("L40").Value = 1333.3334
$ws.Range("M40").Value = -935.7141999999999
$ws.Range("N40").Value = -1683.3334

$ws.Range("H131").Value = 1686.25
$ws.Range("I131").Value = 1498
$ws.Range("J131").Value = 2000
$ws.Range("K131").Value = 4494
$ws.Range("L131").Value = 6000
$ws.Range("M131").Value = 546
$ws.Range("N131").Value = -16080

$ws.Range("H132").Value = 3232
$ws.Range("I132").Value = 1665
$ws.Range("J132").Value = 9500
$ws.Range("K132").Value = 4995
$ws.Range("L132").Value = 28500
$ws.Range("M132").Value = -2465
$ws.Range("N132").Value = -33560

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = ""

$ws.Range("H134").Value = 50000
$ws.Range("J134").Value = 50000
$ws.Range("L134").Value = 50000
$ws.Range("N134").Value = -60140

$ws.Range("H135").Value = 1415.1459
$ws.Range("I135").Value = 946.6
$ws.Range("J135").Value = 8443.333000000001
$ws.Range("K135").Value = 8519.4
$ws.Range("L135").Value = 75989.997
$ws.Range("M135").Value = -5984.4
$ws.Range("N135").Value = -81059.997

$ws.Range("H136").Value = 20000
$ws.Range("J136").Value = 20000
$ws.Range("L136").Value = 20000
$ws.Range("N136").Value = -30200

$ws.Range("H137").Value = 673.4545000000001
$ws.Range("I137").Value = 548.55316
$ws.Range("J137").Value = 982.4211
$ws.Range("K137").Value = 1645.65948
$ws.Range("L137").Value = 2947.2633
$ws.Range("M137").Value = 904.3405199999997
$ws.Range("N137").Value = -8047.263300000001

$ws.Range("H138").Value = 3687.1206
$ws.Range("I138").Value = 2183.4119
$ws.Range("J138").Value = 5817.375
$ws.Range("K138").Value = 6550.2357
$ws.Range("L138").Value = 17452.125
$ws.Range("M138").Value = -1410.2357
$ws.Range("N138").Value = -27732.125

$ws.Range("H139").Value = 78000
$ws.Range("J139").Value = 78000
$ws.Range("L139").Value = 78000
$ws.Range("N139").Value = -88280

$ws.Range("H140").Value = 88950
$ws.Range("J140").Value = 88950
$ws.Range("L140").Value = 88950
$ws.Range("N140").Value = -99310

$ws.Range("H141").Value = 1673.1041
$ws.Range("I141").Value = 1595.5682
$ws.Range("J141").Value = 2526
$ws.Range("K141").Value = 4786.7046
$ws.Range("L141").Value = 7578
$ws.Range("M141").Value = 393.2954
$ws.Range("N141").Value = -17938

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 479859.6
$ws.Range("I32").Value = 3329.9456
$ws.Range("K32").Value = 3329.9456
$ws.Range("M32").Value = -3042.9456

$ws.Range("H74").Value = 1027.8857
$ws.Range("I74").Value = 1182
$ws.Range("J74").Value = 582.6667
$ws.Range("K74").Value = 1182
$ws.Range("L74").Value = 582.6667
$ws.Range("M74").Value = -308
$ws.Range("N74").Value = -2330.6667

$ws.Range("H77").Value = 1027.8857
$ws.Range("I77").Value = 1182
$ws.Range("J77").Value = 582.6667
$ws.Range("K77").Value = 5910
$ws.Range("L77").Value = 2913.3335
$ws.Range("M77").Value = -1542
$ws.Range("N77").Value = -11649.3335

$ws.Range("H122").Value = 2201.5
$ws.Range("I122").Value = 2227.85
$ws.Range("J122").Value = 1938
$ws.Range("K122").Value = 6683.549999999999
$ws.Range("L122").Value = 5814
$ws.Range("M122").Value = -4233.549999999999
$ws.Range("N122").Value = -10714

$ws.Range("H132").Value = 1763.8636
$ws.Range("I132").Value = 1426.1052
$ws.Range("J132").Value = 3903
$ws.Range("K132").Value = 4278.3156
$ws.Range("L132").Value = 11709
$ws.Range("M132").Value = -1748.3156
$ws.Range("N132").Value = -16769

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5523.7
$ws.Range("I134").Value = 1565.1852
$ws.Range("K134").Value = 4695.5556
$ws.Range("M134").Value = -2160.5556

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1426.4667
$ws.Range("I16").Value = 1425.8182
$ws.Range("J16").Value = 1428.25
$ws.Range("K16").Value = 1425.8182
$ws.Range("L16").Value = 1428.25
$ws.Range("M16").Value = -1138.8182
$ws.Range("N16").Value = -2002.25

$ws.Range("H31").Value = 3365.7322
$ws.Range("I31").Value = 1740.4186
$ws.Range("J31").Value = 8741.77
$ws.Range("K31").Value = 1740.4186
$ws.Range("L31").Value = 8741.77
$ws.Range("M31").Value = -1445.4186
$ws.Range("N31").Value = -9331.77

$ws.Range("H34").Value = 3365.7322
$ws.Range("I34").Value = 1740.4186
$ws.Range("J34").Value = 8741.77
$ws.Range("K34").Value = 1740.4186
$ws.Range("L34").Value = 8741.77
$ws.Range("M34").Value = -1538.4186
$ws.Range("N34").Value = -9145.77

$ws.Range("H58").Value = 1309.7567
$ws.Range("I58").Value = 1076.7084
$ws.Range("J58").Value = 1740
$ws.Range("K58").Value = 1076.7084
$ws.Range("L58").Value = 1740
$ws.Range("M58").Value = -873.7084
$ws.Range("N58").Value = -2146

$ws.Range("H94").Value = 5890.4614
$ws.Range("I94").Value = 2050
$ws.Range("J94").Value = 7042.6
$ws.Range("K94").Value = 2050
$ws.Range("L94").Value = 7042.6
$ws.Range("M94").Value = -1599
$ws.Range("N94").Value = -7944.6

$ws.Range("H99").Value = 1484
$ws.Range("I99").Value = 1480
$ws.Range("J99").Value = 1500
$ws.Range("K99").Value = 1480
$ws.Range("L99").Value = 1500
$ws.Range("M99").Value = 18
$ws.Range("N99").Value = -4496

$ws.Range("H105").Value = 1655.8235
$ws.Range("I105").Value = 1124.4546
$ws.Range("K105").Value = 1124.4546
$ws.Range("M105").Value = 622.5454

$ws.Range("H107").Value = 480.8125
$ws.Range("I107").Value = 412.86667
$ws.Range("J107").Value = 1500
$ws.Range("K107").Value = 412.86667
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = 1507.13333
$ws.Range("N107").Value = -5340

$ws.Range("H113").Value = 1426.4667
$ws.Range("I113").Value = 1425.8182
$ws.Range("J113").Value = 1428.25
$ws.Range("K113").Value = 1425.8182
$ws.Range("L113").Value = 1428.25
$ws.Range("M113").Value = 744.1818000000001
$ws.Range("N113").Value = -5768.25

$ws.Range("H126").Value = 1484
$ws.Range("I126").Value = 1480
$ws.Range("J126").Value = 1500
$ws.Range("K126").Value = 4440
$ws.Range("L126").Value = 4500
$ws.Range("M126").Value = -1970
$ws.Range("N126").Value = -9440

$ws.Range("H132").Value = 2466.3684
$ws.Range("I132").Value = 2207.8333
$ws.Range("J132").Value = 2909.5715
$ws.Range("K132").Value = 6623.499899999999
$ws.Range("L132").Value = 8728.7145
$ws.Range("M132").Value = -4093.499899999999
$ws.Range("N132").Value = -13788.7145

$ws.Range("H134").Value = 1948.4166
$ws.Range("I134").Value = 2002.2858
$ws.Range("J134").Value = 1571.3334
$ws.Range("K134").Value = 6006.857400000001
$ws.Range("L134").Value = 4714.0002
$ws.Range("M134").Value = -3471.857400000001
$ws.Range("N134").Value = -9784.0002

$ws.Range("H136").Value = 1309.7567
$ws.Range("I136").Value = 1076.7084
$ws.Range("J136").Value = 1740
$ws.Range("K136").Value = 3230.1252
$ws.Range("L136").Value = 5220
$ws.Range("M136").Value = -680.1251999999999
$ws.Range("N136").Value = -10320

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 37607.37
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 37607.37
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 112822.11
$ws.Range("M121").Value = ""
$ws.Range("N121").Value = -115442.11

$ws.Range("H131").Value = 7247177.5
$ws.Range("J131").Value = 8197566.5
$ws.Range("L131").Value = 24592699.5
$ws.Range("N131").Value = -24602779.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 5000
$ws.Range("J44").Value = 5000
$ws.Range("L44").Value = 5000
$ws.Range("N44").Value = -6192

$ws.Range("H122").Value = 2046.6552
$ws.Range("I122").Value = 1746.2174
$ws.Range("J122").Value = 3198.3333
$ws.Range("K122").Value = 5238.6522
$ws.Range("L122").Value = 9594.999899999999
$ws.Range("M122").Value = -2788.6522
$ws.Range("N122").Value = -14494.9999

$ws.Range("H126").Value = 10418093
$ws.Range("I126").Value = 1281
$ws.Range("J126").Value = 27779446
$ws.Range("K126").Value = 3843
$ws.Range("L126").Value = 83338338
$ws.Range("M126").Value = -1373
$ws.Range("N126").Value = -83343278

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 3500
$ws.Range("J5").Value = 3500
$ws.Range("L5").Value = 3500
$ws.Range("N5").Value = -3726

$ws.Range("H40").Value = 3288.3333
$ws.Range("I40").Value = 3000
$ws.Range("J40").Value = 3432.5
$ws.Range("K40").Value = 3000
$ws.Range("L40").Value = 3432.5
$ws.Range("M40").Value = -2864
$ws.Range("N40").Value = -3704.5

$ws.Range("H132").Value = 2085.4688
$ws.Range("I132").Value = 2160.0408
$ws.Range("J132").Value = 1841.8667
$ws.Range("K132").Value = 6480.1224
$ws.Range("L132").Value = 5525.6001
$ws.Range("M132").Value = -3950.1224
$ws.Range("N132").Value = -10585.6001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 891.2222
$ws.Range("I136").Value = 808.0263
$ws.Range("J136").Value = 1342.8572
$ws.Range("K136").Value = 2424.0789
$ws.Range("L136").Value = 4028.5716
$ws.Range("M136").Value = 125.9211
$ws.Range("N136").Value = -9128.571599999999
